$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the table contents one column to the left (B->A, C->B) and clear column C,
# without disturbing the column widths / sheet dimension metadata.
for ($r = 1; $r -le 6; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $cVal = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 1).Value = $bVal
    $ws.Cells.Item($r, 2).Value = $cVal
    $ws.Cells.Item($r, 3).ClearContents()
}

# Update the selection to reflect the new active cell / selected range
$ws.Range("C1:C6").Select()
